$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J ("2021") is a new year column added after I ("2020"), mirroring
# the same row layout (section-rule rows stay blank, data rows get a value).

# Row 3 - top thin rule under the title; blank cell, style copied from I3.
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial(-4122)

# Row 4 - column headers ("2015" ... "2020"); new header "2021".
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = 2021

# Data rows - plain value assignment inherits the correct style from each
# row's own default format (rows 5/6 are customFormat rows; 8/9/11/12/14-17
# pick up the sheet's general data style automatically).
$ws.Range("J5").Value = 5356.3
$ws.Range("J6").Value = 9.5
$ws.Range("J8").Value = 7.9
$ws.Range("J9").Value = 10.5
$ws.Range("J11").Value = 9.6
$ws.Range("J12").Value = 9.4
$ws.Range("J14").Value = 14.8
$ws.Range("J15").Value = 9.1
$ws.Range("J16").Value = 9.5
$ws.Range("J17").Value = 5.9

# Row 27 - bottom total rule; blank cell, style copied from I27.
$ws.Range("I27").Copy()
$ws.Range("J27").PasteSpecial(-4122)

# Restore the cursor position recorded in the saved file's sheet view.
[void]$ws.Range("L27").Select()
